# GitHub Actions symbol-list refresh (cryptos.xlsx) — re-applies the scraped
# coinranking.com price/volume snapshot onto sheet1:
#   - rows 2-9, 19-50: refreshed "Price" (column D) values, plus two
#     "Volume(1h)" (column E) label tweaks
#   - rows 10-18: a new coin ("One") was inserted ahead of "WazirX", which
#     pushed WazirX..CoinExToken down by one row each (their Coin/Link/
#     Price/Volume columns all shift down one row); the old trailing "One"
#     row is gone because its data now lives at row 10
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text (non-numeric-looking) cell updates
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("E45").Value = "44ACDXExchangeACXT"

# Numeric-looking values that must remain stored as text (matches source t="inlineStr" cells)
$numericTextCells = @{
    "D2" = "250.20"
    "D3" = "23.60"
    "D4" = "5.949"
    "D5" = "0.05928"
    "D6" = "3.429"
    "D7" = "6.574"
    "D8" = "1.329"
    "D9" = "0.7929"
    "D10" = "0.01268"
    "D11" = "0.1483"
    "D12" = "0.07897"
    "D13" = "0.03338"
    "D14" = "0.03039"
    "D15" = "0.09261"
    "D16" = "3.567"
    "D17" = "0.001667"
    "D18" = "0.04784"
    "D19" = "0.006221"
    "D20" = "0.005679"
    "D22" = "0.0001502"
    "D23" = "3.711"
    "D25" = "0.3305"
    "D27" = "0.0006482"
    "D40" = "0.04425"
    "D41" = "0.007036"
    "D42" = "0.1065"
    "D43" = "0.003304"
    "D44" = "0.009391"
    "D45" = "0.002462"
    "D46" = "0.00005901"
    "D47" = "0.00000000751"
    "D48" = "0.9911"
    "D49" = "0.1112"
    "D50" = "0.00002102"
}
foreach ($addr in $numericTextCells.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $numericTextCells[$addr]
}
